# Update column F (dSF) values per the "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 2
    4  = 7
    7  = -1
    8  = 8
    9  = 4
    10 = 3
    11 = -3
    12 = -1
    13 = 1
    15 = -1
    17 = -1
    18 = -1
    19 = 2
    20 = -4
    21 = 1
    22 = -1
    23 = -1
    24 = -5
    25 = 1
    26 = -1
    27 = 1
    28 = 8
    29 = 2
    30 = 3
    32 = -1
    33 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
